# ----------------------------------------------------------------------------
# BIOC_YR_FIN.xlsx - "Doing Updates for Financials"
#
# A new yearly column of financial data is inserted at column D (the most
# recent period). Existing data in D:I shifts right into E:J, the former
# column J (the oldest period, now out of the 7-year window) is dropped, and
# a handful of cells that no longer have a value are marked "NA". This
# mirrors the row-level edits recorded for the BIOC worksheet.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIOC")

# Row 8 - Total Revenue
$ws.Range("D8").Value = 3300
$ws.Range("E8").Value = 5100
$ws.Range("F8").Value = 3200
$ws.Range("G8").Value = 600
$ws.Range("J8").Value = 100

# Row 9 - Cost of Revenue
$ws.Range("D9").Value = 10100
$ws.Range("E9").Value = 9300
$ws.Range("F9").Value = 6900
$ws.Range("G9").Value = 4600
$ws.Range("H9").Value = 2200
$ws.Range("I9").Value = 2300
$ws.Range("J9").Value = 1200

# Row 10 - Gross Profit
$ws.Range("D10").Value = -6800
$ws.Range("E10").Value = -4300
$ws.Range("F10").Value = -3700
$ws.Range("G10").Value = -4000
$ws.Range("H10").Value = -2000
$ws.Range("I10").Value = -2200
$ws.Range("J10").Value = -1100

# Row 12 - Research Development
$ws.Range("D12").Value = 4500
$ws.Range("E12").Value = 3400
$ws.Range("F12").Value = 2700
$ws.Range("G12").Value = 2900
$ws.Range("H12").Value = 4500
$ws.Range("I12").Value = 3100
$ws.Range("J12").Value = 6600

# Row 17 - Total Operating Expenses
$ws.Range("D17").Value = 27500
$ws.Range("E17").Value = 26200
$ws.Range("F17").Value = 21200
$ws.Range("G17").Value = 17000
$ws.Range("H17").Value = 14000
$ws.Range("I17").Value = 8100
$ws.Range("J17").Value = 10600

# Row 18 - Operating Income or Loss
$ws.Range("D18").Value = -24300
$ws.Range("E18").Value = -21200
$ws.Range("F18").Value = -18000
$ws.Range("G18").Value = -16400
$ws.Range("H18").Value = -13900
$ws.Range("I18").Value = -7900
$ws.Range("J18").Value = -10500

# Row 20 - Total Other Income/Expenses Net
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 200
$ws.Range("G20").Value = 100
$ws.Range("H20").Value = -200
$ws.Range("I20").Value = 800
$ws.Range("J20").Value = 400

# Row 21 - Earnings Before Interest And Taxes
$ws.Range("D21").Value = -23700
$ws.Range("E21").Value = -20800
$ws.Range("F21").Value = -17600
$ws.Range("G21").Value = -16100
$ws.Range("H21").Value = -13800
$ws.Range("I21").Value = -6800
$ws.Range("J21").Value = -9800

# Row 22 - Interest Expense
$ws.Range("D22").Value = 300
$ws.Range("F22").Value = 500
$ws.Range("G22").Value = 600
$ws.Range("H22").Value = 1800
$ws.Range("J22").Value = 2100

# Row 23 - Income Before Tax
$ws.Range("D23").Value = -24600
$ws.Range("E23").Value = -21600
$ws.Range("F23").Value = -18400
$ws.Range("G23").Value = -16900
$ws.Range("H23").Value = -15900
$ws.Range("I23").Value = -9200
$ws.Range("J23").Value = -12300

# Row 24 - Income Tax Expense
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = -2600

# Row 26 - Income After Tax
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = -19000
$ws.Range("F26").Value = -18400
$ws.Range("G26").Value = -16900
$ws.Range("H26").Value = -15900
$ws.Range("I26").Value = -9200
$ws.Range("J26").Value = -12300

# Row 27 - Net Income From Continuing Ops
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = -19000
$ws.Range("F27").Value = -18400
$ws.Range("G27").Value = -16900
$ws.Range("H27").Value = -15900
$ws.Range("I27").Value = -9200
$ws.Range("J27").Value = -12300

# Row 29 - Discontinued Operations
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = -2600

# Row 32 - Other Items
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = -100
$ws.Range("F32").Value = -200
$ws.Range("G32").Value = -100
$ws.Range("H32").Value = 200
$ws.Range("I32").Value = -800
$ws.Range("J32").Value = -400

# Row 33 - Net Income
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = -21600
$ws.Range("F33").Value = -18400
$ws.Range("G33").Value = -16900
$ws.Range("H33").Value = -15900
$ws.Range("I33").Value = -9200
$ws.Range("J33").Value = -12300

# Row 35 - Net Income Applicable To Common Shares
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = -21600
$ws.Range("F35").Value = -18400
$ws.Range("G35").Value = -16900
$ws.Range("H35").Value = -15900
$ws.Range("I35").Value = -9200
$ws.Range("J35").Value = -12300

# Row 81 - Net Income (Cash Flow Statement)
$ws.Range("D81").Value = "NA"
$ws.Range("E81").Value = -21600
$ws.Range("F81").Value = -18400
$ws.Range("G81").Value = -16900
$ws.Range("H81").Value = -15900
$ws.Range("I81").Value = -9200
$ws.Range("J81").Value = -12300

# Row 83 - Depreciation
$ws.Range("J83").Value = "NA"

# Row 94 - Total Cash Flows From Investing Activities
$ws.Range("J94").Value = "NA"

# Row 100 - Total Cash Flows From Financing Activities
$ws.Range("J100").Value = "NA"

# Row 101 - Effect Of Exchange Rate Changes
$ws.Range("J101").Value = "NA"
